$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2226.6667
$ws.Range("J43").Value = 2299.3333
$ws.Range("L43").Value = 2299.3333
$ws.Range("N43").Value = -2437.3333

$ws.Range("H76").Value = 2466.3333
$ws.Range("I76").Value = 2466.3333
$ws.Range("K76").Value = 2466.3333
$ws.Range("M76").Value = -2151.3333

$ws.Range("H79").Value = 2466.3333
$ws.Range("I79").Value = 2466.3333
$ws.Range("K79").Value = 2466.3333
$ws.Range("M79").Value = -1374.3333

$ws.Range("H111").Value = 1289.8667
$ws.Range("J111").Value = 1838.5555
$ws.Range("L111").Value = 5515.666499999999
$ws.Range("N111").Value = -11649.6665

$ws.Range("H113").Value = 5810
$ws.Range("I113").Value = 5810
$ws.Range("K113").Value = 5810
$ws.Range("M113").Value = -2556

$ws.Range("H116").Value = 4709.923
$ws.Range("I116").Value = 4712
$ws.Range("J116").Value = 4698.5
$ws.Range("K116").Value = 4712
$ws.Range("L116").Value = 4698.5
$ws.Range("M116").Value = -1270
$ws.Range("N116").Value = -11582.5

$ws.Range("H139").Value = 100780
$ws.Range("J139").Value = 100780
$ws.Range("L139").Value = 100780
$ws.Range("N139").Value = -111060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1118.3
$ws.Range("I97").Value = 1110.1177
$ws.Range("K97").Value = 1110.1177
$ws.Range("M97").Value = -614.1177

$ws.Range("H110").Value = 1913.625
$ws.Range("I110").Value = 1756.5714
$ws.Range("K110").Value = 1756.5714
$ws.Range("M110").Value = 288.4286

$ws.Range("H122").Value = 1997.5
$ws.Range("I122").Value = 1997.5
$ws.Range("K122").Value = 5992.5
$ws.Range("M122").Value = -3542.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1037.4615
$ws.Range("I94").Value = 808
$ws.Range("K94").Value = 808
$ws.Range("M94").Value = -357

$ws.Range("H107").Value = 6221.8184
$ws.Range("I107").Value = 6045
$ws.Range("K107").Value = 6045
$ws.Range("M107").Value = -4125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2420.5
$ws.Range("I62").Value = 2302
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 2302
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -1678
$ws.Range("N62").Value = -4498

$ws.Range("H65").Value = 2420.5
$ws.Range("I65").Value = 2302
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 11510
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -8390
$ws.Range("N65").Value = -22490

$ws.Range("H105").Value = 3598.5
$ws.Range("I105").Value = 3598.5
$ws.Range("K105").Value = 3598.5
$ws.Range("M105").Value = -1851.5

$ws.Range("H122").Value = 1109.1111
$ws.Range("I122").Value = 1109.1111
$ws.Range("K122").Value = 3327.3333
$ws.Range("M122").Value = -877.3333000000002

$ws.Range("H132").Value = 2161
$ws.Range("I132").Value = 2161
$ws.Range("K132").Value = 6483
$ws.Range("M132").Value = -3953

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1250.75
$ws.Range("I5").Value = 1341.6
$ws.Range("J5").Value = 1185.8572
$ws.Range("K5").Value = 4024.8
$ws.Range("L5").Value = 3557.5716
$ws.Range("M5").Value = -3912.8
$ws.Range("N5").Value = -3781.5716

$ws.Range("H7").Value = 1033.3334
$ws.Range("J7").Value = 1050
$ws.Range("L7").Value = 3150
$ws.Range("N7").Value = -3374

$ws.Range("H131").Value = 1712.8
$ws.Range("J131").Value = 2995.6667
$ws.Range("L131").Value = 8987.000100000001
$ws.Range("N131").Value = -19067.0001

$ws.Range("H135").Value = 1250.75
$ws.Range("I135").Value = 1341.6
$ws.Range("J135").Value = 1185.8572
$ws.Range("K135").Value = 12074.4
$ws.Range("L135").Value = 10672.7148
$ws.Range("M135").Value = -9539.4
$ws.Range("N135").Value = -15742.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 731.6667
$ws.Range("I80").Value = 531.1667
$ws.Range("K80").Value = 531.1667
$ws.Range("M80").Value = 466.8333

$ws.Range("H83").Value = 731.6667
$ws.Range("I83").Value = 531.1667
$ws.Range("K83").Value = 2655.8335
$ws.Range("M83").Value = 2336.1665

$ws.Range("H97").Value = 998.2857
$ws.Range("I97").Value = 699
$ws.Range("K97").Value = 699
$ws.Range("M97").Value = -203

$ws.Range("H102").Value = 2718.353
$ws.Range("I102").Value = 2414.8
$ws.Range("J102").Value = 4995
$ws.Range("K102").Value = 2414.8
$ws.Range("L102").Value = 4995
$ws.Range("M102").Value = -792.8000000000002
$ws.Range("N102").Value = -8239

$ws.Range("H107").Value = 374.25
$ws.Range("I107").Value = 199
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 199
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1721
$ws.Range("N107").Value = -4740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1295.5454
$ws.Range("I16").Value = 1059
$ws.Range("J16").Value = 1926.3334
$ws.Range("K16").Value = 1059
$ws.Range("L16").Value = 1926.3334
$ws.Range("M16").Value = -889
$ws.Range("N16").Value = -2266.3334

$ws.Range("H68").Value = 1333
$ws.Range("I68").Value = 1333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1333
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -584
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1333
$ws.Range("I71").Value = 1333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6665
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2921
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 5964.143
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 6791.5
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 6791.5
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -7513.5

$ws.Range("H85").Value = 5964.143
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 6791.5
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 6791.5
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -9287.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 839.4286
$ws.Range("I107").Value = 829.4
$ws.Range("J107").Value = 864.5
$ws.Range("K107").Value = 2488.2
$ws.Range("L107").Value = 2593.5
$ws.Range("M107").Value = -568.1999999999998
$ws.Range("N107").Value = -6433.5

$ws.Range("H113").Value = 713.9231
$ws.Range("J113").Value = 793.875
$ws.Range("L113").Value = 2381.625
$ws.Range("N113").Value = -6721.625

$ws.Range("H122").Value = 1373.4117
$ws.Range("I122").Value = 1234.1428
$ws.Range("J122").Value = 1470.9
$ws.Range("K122").Value = 3702.4284
$ws.Range("L122").Value = 4412.700000000001
$ws.Range("M122").Value = -1252.4284
$ws.Range("N122").Value = -9312.700000000001

$ws.Range("H126").Value = 7999

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 3288.2222
$ws.Range("I132").Value = 3156.5715
$ws.Range("K132").Value = 9469.7145
$ws.Range("M132").Value = -6939.7145
